$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 10291.667
$ws.Cells.Item(40, 9).Value = 6747.3335
$ws.Cells.Item(40, 11).Value = 6747.3335
$ws.Cells.Item(40, 13).Value = -6572.3335
$ws.Cells.Item(82, 8).Value = 250.66667
$ws.Cells.Item(82, 9).Value = 151
$ws.Cells.Item(82, 10).Value = 450
$ws.Cells.Item(82, 11).Value = 453
$ws.Cells.Item(82, 12).Value = 1350
$ws.Cells.Item(82, 13).Value = -47
$ws.Cells.Item(82, 14).Value = -2162
$ws.Cells.Item(85, 8).Value = 250.66667
$ws.Cells.Item(85, 9).Value = 151
$ws.Cells.Item(85, 10).Value = 450
$ws.Cells.Item(85, 11).Value = 453
$ws.Cells.Item(85, 12).Value = 1350
$ws.Cells.Item(85, 13).Value = 951
$ws.Cells.Item(85, 14).Value = -4158
$ws.Cells.Item(132, 8).Value = 2469.1892
$ws.Cells.Item(132, 9).Value = 2656.6667
$ws.Cells.Item(132, 10).Value = 922.5
$ws.Cells.Item(132, 11).Value = 7970.000100000001
$ws.Cells.Item(132, 12).Value = 2767.5
$ws.Cells.Item(132, 13).Value = -5440.000100000001
$ws.Cells.Item(132, 14).Value = -7827.5
$ws.Cells.Item(137, 8).Value = 3280.805
$ws.Cells.Item(137, 9).Value = 3008.5557
$ws.Cells.Item(137, 10).Value = 3357.375
$ws.Cells.Item(137, 11).Value = 9025.667099999999
$ws.Cells.Item(137, 12).Value = 10072.125
$ws.Cells.Item(137, 13).Value = -6475.667099999999
$ws.Cells.Item(137, 14).Value = -15172.125
$ws.Cells.Item(138, 8).Value = 3418.75
$ws.Cells.Item(138, 9).Value = 2631.7646
$ws.Cells.Item(138, 10).Value = 4635
$ws.Cells.Item(138, 11).Value = 7895.293799999999
$ws.Cells.Item(138, 12).Value = 13905
$ws.Cells.Item(138, 13).Value = -2755.293799999999
$ws.Cells.Item(138, 14).Value = -24185

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 6946936.5
$ws.Cells.Item(74, 9).Value = 8132729
$ws.Cells.Item(74, 10).Value = 1581.1428
$ws.Cells.Item(74, 11).Value = 8132729
$ws.Cells.Item(74, 12).Value = 1581.1428
$ws.Cells.Item(74, 13).Value = -8131855
$ws.Cells.Item(74, 14).Value = -3329.1428
$ws.Cells.Item(77, 8).Value = 6946936.5
$ws.Cells.Item(77, 9).Value = 8132729
$ws.Cells.Item(77, 10).Value = 1581.1428
$ws.Cells.Item(77, 11).Value = 40663645
$ws.Cells.Item(77, 12).Value = 7905.714
$ws.Cells.Item(77, 13).Value = -40659277
$ws.Cells.Item(77, 14).Value = -16641.714
$ws.Cells.Item(135, 8).Value = 54999.285
$ws.Cells.Item(135, 10).Value = 54999.285
$ws.Cells.Item(135, 12).Value = 54999.285
$ws.Cells.Item(135, 14).Value = -65139.285
$ws.Cells.Item(140, 8).Value = 88550.44500000001
$ws.Cells.Item(140, 10).Value = 88550.44500000001
$ws.Cells.Item(140, 12).Value = 88550.44500000001
$ws.Cells.Item(140, 14).Value = -98910.44500000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2957.1428
$ws.Cells.Item(99, 9).Value = 2616.6667
$ws.Cells.Item(99, 10).Value = 5000
$ws.Cells.Item(99, 11).Value = 2616.6667
$ws.Cells.Item(99, 12).Value = 5000
$ws.Cells.Item(99, 13).Value = -1118.6667
$ws.Cells.Item(99, 14).Value = -7996
$ws.Cells.Item(107, 8).Value = 3550.3076
$ws.Cells.Item(107, 9).Value = 2105.875
$ws.Cells.Item(107, 10).Value = 5861.4
$ws.Cells.Item(107, 11).Value = 2105.875
$ws.Cells.Item(107, 12).Value = 5861.4
$ws.Cells.Item(107, 13).Value = -185.875
$ws.Cells.Item(107, 14).Value = -9701.4
$ws.Cells.Item(134, 8).Value = 4045.7878
$ws.Cells.Item(134, 9).Value = 2061.75
$ws.Cells.Item(134, 11).Value = 6185.25
$ws.Cells.Item(134, 13).Value = -3650.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 22985.223
$ws.Cells.Item(31, 9).Value = 3580.75
$ws.Cells.Item(31, 11).Value = 3580.75
$ws.Cells.Item(31, 13).Value = -3285.75
$ws.Cells.Item(34, 8).Value = 22985.223
$ws.Cells.Item(34, 9).Value = 3580.75
$ws.Cells.Item(34, 11).Value = 3580.75
$ws.Cells.Item(34, 13).Value = -3378.75
$ws.Cells.Item(62, 8).Value = 13210.777
$ws.Cells.Item(62, 9).Value = 14222
$ws.Cells.Item(62, 10).Value = 12401.8
$ws.Cells.Item(62, 11).Value = 14222
$ws.Cells.Item(62, 12).Value = 12401.8
$ws.Cells.Item(62, 13).Value = -13598
$ws.Cells.Item(62, 14).Value = -13649.8
$ws.Cells.Item(65, 8).Value = 13210.777
$ws.Cells.Item(65, 9).Value = 14222
$ws.Cells.Item(65, 10).Value = 12401.8
$ws.Cells.Item(65, 11).Value = 71110
$ws.Cells.Item(65, 12).Value = 62009
$ws.Cells.Item(65, 13).Value = -67990
$ws.Cells.Item(65, 14).Value = -68249
$ws.Cells.Item(125, 8).Value = 61605.7
$ws.Cells.Item(125, 10).Value = 61605.7
$ws.Cells.Item(125, 12).Value = 61605.7
$ws.Cells.Item(125, 14).Value = -66525.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 3089.2307
$ws.Cells.Item(68, 10).Value = 3152.8
$ws.Cells.Item(68, 12).Value = 9458.400000000001
$ws.Cells.Item(68, 14).Value = -11080.4
$ws.Cells.Item(71, 8).Value = 3089.2307
$ws.Cells.Item(71, 10).Value = 3152.8
$ws.Cells.Item(71, 12).Value = 28375.2
$ws.Cells.Item(71, 14).Value = -36487.2
$ws.Cells.Item(88, 8).Value = 16338.667
$ws.Cells.Item(88, 10).Value = 17008
$ws.Cells.Item(88, 12).Value = 51024
$ws.Cells.Item(88, 14).Value = -51880
$ws.Cells.Item(91, 8).Value = 16338.667
$ws.Cells.Item(91, 10).Value = 17008
$ws.Cells.Item(91, 12).Value = 51024
$ws.Cells.Item(91, 14).Value = -53988
$ws.Cells.Item(107, 8).Value = 1950.3914
$ws.Cells.Item(107, 10).Value = 2197.611
$ws.Cells.Item(107, 12).Value = 6592.833
$ws.Cells.Item(107, 14).Value = -10432.833
$ws.Cells.Item(139, 8).Value = 3796.7585
$ws.Cells.Item(139, 9).Value = 1726.35
$ws.Cells.Item(139, 10).Value = 8397.666999999999
$ws.Cells.Item(139, 11).Value = 5179.049999999999
$ws.Cells.Item(139, 12).Value = 25193.001
$ws.Cells.Item(139, 13).Value = -39.04999999999927
$ws.Cells.Item(139, 14).Value = -35473.001
$ws.Cells.Item(140, 8).Value = 852.625
$ws.Cells.Item(140, 9).Value = 852.625
$ws.Cells.Item(140, 11).Value = 2557.875
$ws.Cells.Item(140, 13).Value = 2622.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6515.125
$ws.Cells.Item(70, 9).Value = 6416.75
$ws.Cells.Item(70, 11).Value = 6416.75
$ws.Cells.Item(70, 13).Value = -6146.75
$ws.Cells.Item(73, 8).Value = 6515.125
$ws.Cells.Item(73, 9).Value = 6416.75
$ws.Cells.Item(73, 11).Value = 6416.75
$ws.Cells.Item(73, 13).Value = -5480.75
$ws.Cells.Item(80, 8).Value = 12533.667
$ws.Cells.Item(80, 9).Value = 4949.5
$ws.Cells.Item(80, 10).Value = 18601
$ws.Cells.Item(80, 11).Value = 4949.5
$ws.Cells.Item(80, 12).Value = 18601
$ws.Cells.Item(80, 13).Value = -3951.5
$ws.Cells.Item(80, 14).Value = -20597
$ws.Cells.Item(83, 8).Value = 12533.667
$ws.Cells.Item(83, 9).Value = 4949.5
$ws.Cells.Item(83, 10).Value = 18601
$ws.Cells.Item(83, 11).Value = 24747.5
$ws.Cells.Item(83, 12).Value = 93005
$ws.Cells.Item(83, 13).Value = -19755.5
$ws.Cells.Item(83, 14).Value = -102989
$ws.Cells.Item(122, 8).Value = 8771.467000000001
$ws.Cells.Item(122, 9).Value = 4514.3335
$ws.Cells.Item(122, 11).Value = 13543.0005
$ws.Cells.Item(122, 13).Value = -11093.0005
$ws.Cells.Item(132, 8).Value = 6823.607
$ws.Cells.Item(132, 9).Value = 1378.1333
$ws.Cells.Item(132, 10).Value = 13106.846
$ws.Cells.Item(132, 11).Value = 4134.3999
$ws.Cells.Item(132, 12).Value = 39320.538
$ws.Cells.Item(132, 13).Value = -1604.3999
$ws.Cells.Item(132, 14).Value = -44380.538

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 650.25
$ws.Cells.Item(35, 9).Value = 650.25
$ws.Cells.Item(35, 11).Value = 650.25
$ws.Cells.Item(35, 13).Value = -314.25
$ws.Cells.Item(40, 8).Value = 7843.25
$ws.Cells.Item(40, 9).Value = 5348.9
$ws.Cells.Item(40, 10).Value = 12000.5
$ws.Cells.Item(40, 11).Value = 5348.9
$ws.Cells.Item(40, 12).Value = 12000.5
$ws.Cells.Item(40, 13).Value = -5212.9
$ws.Cells.Item(40, 14).Value = -12272.5
$ws.Cells.Item(68, 8).Value = 4614.615
$ws.Cells.Item(68, 9).Value = 2817
$ws.Cells.Item(68, 10).Value = 14501.5
$ws.Cells.Item(68, 11).Value = 2817
$ws.Cells.Item(68, 12).Value = 14501.5
$ws.Cells.Item(68, 13).Value = -2068
$ws.Cells.Item(68, 14).Value = -15999.5
$ws.Cells.Item(71, 8).Value = 4614.615
$ws.Cells.Item(71, 9).Value = 2817
$ws.Cells.Item(71, 10).Value = 14501.5
$ws.Cells.Item(71, 11).Value = 14085
$ws.Cells.Item(71, 12).Value = 72507.5
$ws.Cells.Item(71, 13).Value = -10341
$ws.Cells.Item(71, 14).Value = -79995.5
$ws.Cells.Item(100, 8).Value = 7500
$ws.Cells.Item(100, 9).Value = 2998.5
$ws.Cells.Item(100, 11).Value = 2998.5
$ws.Cells.Item(100, 13).Value = -2457.5
$ws.Cells.Item(132, 8).Value = 5648.2856
$ws.Cells.Item(132, 9).Value = 4440.6665
$ws.Cells.Item(132, 10).Value = 7258.4443
$ws.Cells.Item(132, 11).Value = 13321.9995
$ws.Cells.Item(132, 12).Value = 21775.3329
$ws.Cells.Item(132, 13).Value = -10791.9995
$ws.Cells.Item(132, 14).Value = -26835.3329
$ws.Cells.Item(136, 8).Value = 5419.656
$ws.Cells.Item(136, 9).Value = 6455.2856
$ws.Cells.Item(136, 10).Value = 4540.9395
$ws.Cells.Item(136, 11).Value = 19365.8568
$ws.Cells.Item(136, 12).Value = 13622.8185
$ws.Cells.Item(136, 13).Value = -16815.8568
$ws.Cells.Item(136, 14).Value = -18722.8185

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 5969.222
$ws.Cells.Item(74, 10).Value = 5611.1665
$ws.Cells.Item(74, 12).Value = 5611.1665
$ws.Cells.Item(74, 14).Value = -7483.1665
$ws.Cells.Item(77, 8).Value = 5969.222
$ws.Cells.Item(77, 10).Value = 5611.1665
$ws.Cells.Item(77, 12).Value = 16833.4995
$ws.Cells.Item(77, 14).Value = -26193.4995
$ws.Cells.Item(122, 8).Value = 4141.96
$ws.Cells.Item(122, 9).Value = 2055.1052
$ws.Cells.Item(122, 10).Value = 10750.333
$ws.Cells.Item(122, 11).Value = 6165.3156
$ws.Cells.Item(122, 12).Value = 32250.999
$ws.Cells.Item(122, 13).Value = -3715.3156
$ws.Cells.Item(122, 14).Value = -37150.999
$ws.Cells.Item(126, 8).Value = 7487.25
$ws.Cells.Item(126, 10).Value = 7487.25
$ws.Cells.Item(126, 12).Value = 22461.75
$ws.Cells.Item(126, 14).Value = -27401.75
$ws.Cells.Item(132, 8).Value = 10578.556
$ws.Cells.Item(132, 9).Value = 7660.6
$ws.Cells.Item(132, 10).Value = 14226
$ws.Cells.Item(132, 11).Value = 22981.8
$ws.Cells.Item(132, 12).Value = 42678
$ws.Cells.Item(132, 13).Value = -20451.8
$ws.Cells.Item(132, 14).Value = -47738
$ws.Cells.Item(135, 8).Value = 108336.73
$ws.Cells.Item(135, 10).Value = 108336.73
$ws.Cells.Item(135, 12).Value = 108336.73
$ws.Cells.Item(135, 14).Value = -118476.73
$ws.Cells.Item(136, 8).Value = 3825.5
$ws.Cells.Item(136, 9).Value = 2144.5938
$ws.Cells.Item(136, 11).Value = 6433.7814
$ws.Cells.Item(136, 13).Value = -3883.7814
